# Fix/normalize the ordering of entries in the "Recorded By" (column G) cells
# on the "Session Analysis Results" sheet. This reorders specific comma
# separated lists of recorder identities to match the canonical ordering
# used upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Ordered, case-sensitive list of exact "before" text -> "after" text.
# (A plain hashtable does case-insensitive key lookups by default, which
# would be unsafe here since "System" vs "system" are meaningfully
# different tokens in this data.)
$replacements = New-Object System.Collections.Generic.List[object]
[void]$replacements.Add(@("system, System, backup@backdoor.com", "System, backup@backdoor.com, system"))
[void]$replacements.Add(@("dnasr281@gmail.com, System", "System, dnasr281@gmail.com"))
[void]$replacements.Add(@("dnasr281@gmail.com, admin@admin.com", "admin@admin.com, dnasr281@gmail.com"))

# Determine the used range of the sheet so we know how many rows to scan.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Column G holds the "Recorded By" values; data starts at row 2 (row 1 is the header).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value2

    if ($null -ne $value) {
        foreach ($pair in $replacements) {
            if ($value.Equals($pair[0])) {
                $cell.Value2 = $pair[1]
                break
            }
        }
    }
}
